$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.998769987699877
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0.995079950799508
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("B6").Value = 0.9987684729064039
$ws.Range("B7").Value = 0.9987684729064039
$ws.Range("C7").Value = 0.9975369458128078
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.9938423645320197
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 0.9963054187192119
$ws.Range("C10").Value = 1
$ws.Range("B11").Value = 0.9987684729064039
$ws.Range("C11").Value = 1
$ws.Range("A12").Value = 0.860048640282998
$ws.Range("B12").Value = 0.8642493919964626
$ws.Range("C12").Value = 0.862701746628344
$ws.Range("D12").Value = 0.859606455892107
$ws.Range("A13").Value = 0.8671235905372541
$ws.Range("B13").Value = 0.8691134202962636
$ws.Range("C13").Value = 0.8673446827326995
$ws.Range("D13").Value = 0.8655759451691355
$ws.Range("A14").Value = 0.8639982308712959
$ws.Range("B14").Value = 0.8639982308712959
$ws.Range("C14").Value = 0.8659885006634233
$ws.Range("D14").Value = 0.8648827952233525
$ws.Range("A15").Value = 0.8701901813356921
$ws.Range("B15").Value = 0.8739495798319328
$ws.Range("C15").Value = 0.8735072976559045
$ws.Range("D15").Value = 0.8708536045997346
$ws.Range("A16").Value = 0.8586908447589562
$ws.Range("B16").Value = 0.8624502432551968
$ws.Range("C16").Value = 0.8644405130473242
$ws.Range("D16").Value = 0.8567005749668288
$ws.Range("A17").Value = 0.854046881910659
$ws.Range("B17").Value = 0.8604599734630695
$ws.Range("C17").Value = 0.8620079610791685
$ws.Range("D17").Value = 0.8564794338788146
$ws.Range("A18").Value = 0.8693056169836355
$ws.Range("B18").Value = 0.8717381689517912
$ws.Range("C18").Value = 0.8726227333038479
$ws.Range("D18").Value = 0.8679787704555506
$ws.Range("A19").Value = 0.8646616541353384
$ws.Range("B19").Value = 0.871517027863777
$ws.Range("C19").Value = 0.871517027863777
$ws.Range("D19").Value = 0.8657673595754091
$ws.Range("A20").Value = 0.871295886775763
$ws.Range("B20").Value = 0.8759398496240601
$ws.Range("C20").Value = 0.8757187085360459
$ws.Range("D20").Value = 0.8704113224237063
$ws.Range("A21").Value = 0.868421052631579
$ws.Range("B21").Value = 0.877045555064131
$ws.Range("C21").Value = 0.8754975674480319
$ws.Range("D21").Value = 0.871295886775763
$ws.Range("B22").Value = 0.798581560283688
$ws.Range("C22").Value = 0.8042553191489362
$ws.Range("B23").Value = 0.8156028368794326
$ws.Range("C23").Value = 0.8070921985815603
$ws.Range("B24").Value = 0.7900709219858156
$ws.Range("C24").Value = 0.7843971631205674
$ws.Range("B25").Value = 0.7769886363636364
$ws.Range("B26").Value = 0.7784090909090909
$ws.Range("C26").Value = 0.7769886363636364
$ws.Range("B27").Value = 0.7926136363636364
$ws.Range("C27").Value = 0.7897727272727273
$ws.Range("C28").Value = 0.8025568181818182
$ws.Range("B29").Value = 0.796875
$ws.Range("C29").Value = 0.796875
$ws.Range("B30").Value = 0.7883522727272727
$ws.Range("C30").Value = 0.7926136363636364
$ws.Range("B31").Value = 0.7883522727272727
$ws.Range("A32").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 1
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("A34").Value = 1
$ws.Range("A35").Value = 1
$ws.Range("A37").Value = 1
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("A38").Value = 1
$ws.Range("B38").Value = 1
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("A39").Value = 0.999
$ws.Range("B39").Value = 0.999
$ws.Range("C39").Value = 0.999
$ws.Range("D39").Value = 0.998
$ws.Range("A40").Value = 0.999
$ws.Range("A41").Value = 0.998
$ws.Range("B41").Value = 0.998
$ws.Range("C41").Value = 0.997
$ws.Range("D41").Value = 0.998
